$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, shifting existing rows 161:179 down to 162:180
$ws.Rows(161).Insert()

# Populate the newly inserted row 161 with the new weekly price record
$ws.Range("A161").Value = 11
$ws.Range("B161").Value = "Vega Monumental Concepción"
$ws.Range("C161").Value = "Bíobío"
$ws.Range("D161").Value = 44946
$ws.Range("E161").Value = 8
$ws.Range("F161").Value = 100112043
$ws.Range("G161").Value = "Pepino ensalada"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 270
$ws.Range("K161").Value = 8000
$ws.Range("L161").Value = 9000
$ws.Range("M161").Value = 8556
$ws.Range("N161").Value = "$/caja 60 unidades"
$ws.Range("O161").Value = "Región de Arica y Parinacota"
$ws.Range("P161").Value = 143
$ws.Range("Q161").Value = 60
$ws.Range("R161").Value = "Hortaliza"
